# Replace the calibration-request data rows (A2:G16) with the updated
# equipment list, and clear the rows that no longer have data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (Equipment, Manufacturer, Model, Serial No., ID No., Certificate No.)
# Calibration Date (column F) is left untouched - every row already has 45812.
$data = @(
    @("REFRIGERATOR",           "HAIER",      "HXC-158",          "BE06L8E1T00B2D560005", "KBNH00369", "CH25062727"),
    @("FREEZER",                "HAIER",      "BD-151C",          "B30KE6E1X00BKM9F0069", "KBNH00370", "CH25062728"),
    @("WATER BATH",             "WISEBATH",   "WB-22",            "0400702137M019",       "KBNH00371", "WB25062729"),
    @("CENTRIFUGE",             "HETTICH",    "EBA 21",           "0000146-01-00",        "KBNH00372", "CF25062730"),
    @("CENTRIFUGE",             "HETTICH",    "EBA 20",           "0120977-07",           "KBNH00384", "CF25062731"),
    @("CENTRIFUGE",             "HETTICH",    "EBA 20",           "0120976-07",           "KBNH00388", "CF25062732"),
    @("CENTRIFUGE",             "HETTICH",    "EBA 20",           "0120978-07",           "KBNH00389", "CF25062733"),
    @("HEMATOCRIT CENTRIFUGE",  "HETTICH",    "HAEMATOKRIT 210",  "0013200-03-00",        "KBNH00390", "CF25062734"),
    @("REFRIGERATOR",           "Z-COOL",     "ZCO-2DQ38",        "ZM-2DGB-56049",        "KBNH00392", "CH25062735"),
    @("ROTATOR",                "WISESHAKE",  "SHO-1D",           "00988981349002",       "KBNH00396", "CF25062736")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    # Serial No. column must stay text. Most values already contain letters
    # so Excel keeps them as text automatically; the one purely-numeric-looking
    # value (with leading zeros) needs an explicit text format so it isn't
    # reinterpreted as a number and doesn't lose its leading zeros.
    $serial = $item[3]
    if ($serial -match '^[0-9]+$') {
        $ws.Cells.Item($row, 4).NumberFormat = "@"
    }
    $ws.Cells.Item($row, 4).Value = $serial
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 7).Value = $item[5]
    $row++
}

# Rows 12-16 no longer have equipment entries; clear their contents (keep formatting).
for ($r = 12; $r -le 16; $r++) {
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 7)).ClearContents()
}

# Update the remembered selection to match the saved view.
$ws.Range("B15").Select()
